$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range('D2').Value = '57.249.91'
$ws.Range('E2').Value = '  -2.85%  '

$ws.Range('D3').Value = '2.552.05'
$ws.Range('E3').Value = '  -4.10%  '

$ws.Range('E4').Value = '  -0.01%  '

Set-TextValue $ws.Range('D5') '517.49'
$ws.Range('E5').Value = '  -1.25%  '

Set-TextValue $ws.Range('D6') '141.11'
$ws.Range('E6').Value = '  -2.17%  '

Set-TextValue $ws.Range('D7') '0.998'
$ws.Range('E7').Value = '  -0.10%  '

$ws.Range('E8').Value = '  -2.06%  '

$ws.Range('D9').Value = '2.559.16'
$ws.Range('E9').Value = '  -4.20%  '

$ws.Range('E10').Value = '  -5.74%  '

$ws.Range('E11').Value = '  -3.43%  '

Set-TextValue $ws.Range('D12') '0.323'
$ws.Range('E12').Value = '  -3.67%  '

$ws.Range('E13').Value = '  -0.39%  '

$ws.Range('D14').Value = '3.002.15'
$ws.Range('E14').Value = '  -4.11%  '

$ws.Range('D15').Value = '57.243.07'
$ws.Range('E15').Value = '  -2.85%  '

Set-TextValue $ws.Range('D16') '20.01'
$ws.Range('E16').Value = '  -4.74%  '

$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '2.585.95'
$ws.Range('E17').Value = '  -3.18%  '

$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Range('D18') '0.0000132'
$ws.Range('E18').Value = '  -3.12%  '

Set-TextValue $ws.Range('D19') '330.95'
$ws.Range('E19').Value = '  -2.28%  '

$ws.Range('E20').Value = '  -3.19%  '

$ws.Range('E21').Value = '  -2.78%  '

Set-TextValue $ws.Range('D22') '6.16'
$ws.Range('E22').Value = '  -3.28%  '

$ws.Range('E23').Value = '  -0.03%  '

Set-TextValue $ws.Range('D24') '64.75'
$ws.Range('E24').Value = '  +0.29%  '

$ws.Range('E25').Value = '  +0.89%  '

Set-TextValue $ws.Range('D26') '0.997'
$ws.Range('E26').Value = '  -0.06%  '

$ws.Range('E27').Value = '  -5.10%  '

$ws.Range('D28').Value = '2.657.12'
$ws.Range('E28').Value = '  -4.61%  '

Set-TextValue $ws.Range('D29') '6.90'
$ws.Range('E29').Value = '  -3.43%  '

Set-TextValue $ws.Range('D30') '0.998'
$ws.Range('E30').Value = '  -0.03%  '

$ws.Range('E31').Value = '  -8.15%  '

Set-TextValue $ws.Range('D32') '6.24'
$ws.Range('E32').Value = '  -6.37%  '

$ws.Range('E33').Value = '  -1.72%  '

$ws.Range('B34').Value = 'Monero'
$ws.Range('C34').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range('D34') '148.61'
$ws.Range('E34').Value = '  -1.22%  '

$ws.Range('B35').Value = 'EthereumClassic'
$ws.Range('C35').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range('D35') '18.46'
$ws.Range('E35').Value = '  -2.40%  '

$ws.Range('E36').Value = '  -4.24%  '

$ws.Range('E37').Value = '  -4.83%  '

$ws.Range('E38').Value = '  -7.51%  '

Set-TextValue $ws.Range('D39') '35.61'
$ws.Range('E39').Value = '  -3.29%  '

$ws.Range('E40').Value = '  -5.68%  '

$ws.Range('E41').Value = '  -2.10%  '

$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range('D42') '3.47'
$ws.Range('E42').Value = '  -3.21%  '

$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws.Range('D43') '0.998'
$ws.Range('E43').Value = '  -0.13%  '

Set-TextValue $ws.Range('D44') '10.63'
$ws.Range('E44').Value = '  -0.39%  '

$ws.Range('B45').Value = 'Bittensor'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws.Range('D45') '265.61'
$ws.Range('E45').Value = '  -3.43%  '

$ws.Range('B46').Value = 'Stellar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range('D46') '0.0949'
$ws.Range('E46').Value = '  -2.08%  '

Set-TextValue $ws.Range('D47') '0.579'
$ws.Range('E47').Value = '  -5.86%  '

Set-TextValue $ws.Range('D48') '18.59'
$ws.Range('E48').Value = '  -6.21%  '

Set-TextValue $ws.Range('D49') '0.0513'
$ws.Range('E49').Value = '  -3.76%  '

$ws.Range('D50').Value = '1.958.57'
$ws.Range('E50').Value = '  -4.57%  '

$ws.Range('E51').Value = '  -4.18%  '
